# facilities.xlsx — re-saved from Excel (Mac) after review: the active
# selection moved to D2 and the eight data columns were sized to fit their
# contents (AutoFit-style column widths). No cell values changed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move / record the active selection on Sheet1 (was A1, now D2).
$ws.Range("D2").Select()

# Column widths (A:H), matching the widths Excel wrote after auto-fitting
# the content of each column. The ColumnWidth COM property is expressed in
# "characters"; values below are chosen so the underlying stored column
# width (in the saved XML) lands on the target value.
$ws.Columns.Item(1).ColumnWidth = 2.3333333333333335   # A - id column
$ws.Columns.Item(2).ColumnWidth = 4.333333333333333    # B - price
$ws.Columns.Item(3).ColumnWidth = 54.666666666666664    # C - facility
$ws.Columns.Item(4).ColumnWidth = 51.0                  # D - address
$ws.Columns.Item(5).ColumnWidth = 7.5                   # E - condition
$ws.Columns.Item(6).ColumnWidth = 4.833333333333333     # F - zip
$ws.Columns.Item(7).ColumnWidth = 5.0                   # G - radius
$ws.Columns.Item(8).ColumnWidth = 10.333333333333334    # H - city
